# "stateless entities outside the US"
#
# The income table gains two new "IMF (20%)" columns (Sales / Sales+Emp) for
# each of the three metric blocks (M_%cit: B:I, M_ETR: J:Q, M_PL: R:Y).
# They are inserted right before the existing "IMF - Sales" / "IMF - Sales + Emp"
# columns, while the old "OECD (20%) - Sales" / "OECD (20%) - Sales + Emp"
# columns are dropped. Net effect per 8-column block (columns 3-4-5-6 of each
# block, i.e. D/E/F/G, L/M/N/O, T/U/V/W):
#   col3 -> "IMF (20%) - Sales"        (new data)
#   col4 -> "IMF (20%) - Sales + Emp"  (new data)
#   col5 -> "IMF - Sales"              (reuses what used to be in col3)
#   col6 -> "IMF - Sales + Emp"        (reuses what used to be in col4)
# columns 1,2,7,8 of each block (GFA-Sales, GFA-Sales+Emp, OECD-Sales,
# OECD-Sales+Emp) are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 headers: relabel the 3rd/4th/5th/6th column of each 8-col block ----
$ws.Range("D2").Value = "IMF (20%) - Sales"
$ws.Range("E2").Value = "IMF (20%) - Sales + Emp"
$ws.Range("F2").Value = "IMF - Sales"
$ws.Range("G2").Value = "IMF - Sales + Emp"

$ws.Range("L2").Value = "IMF (20%) - Sales"
$ws.Range("M2").Value = "IMF (20%) - Sales + Emp"
$ws.Range("N2").Value = "IMF - Sales"
$ws.Range("O2").Value = "IMF - Sales + Emp"

$ws.Range("T2").Value = "IMF (20%) - Sales"
$ws.Range("U2").Value = "IMF (20%) - Sales + Emp"
$ws.Range("V2").Value = "IMF - Sales"
$ws.Range("W2").Value = "IMF - Sales + Emp"

# ---- Data rows 4-8 ----

# Row 4 (High Income)
$ws.Range("D4").Value = 1.093635869920273
$ws.Range("E4").Value = 0.9072691627789734
$ws.Range("F4").Value = 5.468179349601362
$ws.Range("G4").Value = 4.536345813894876
$ws.Range("N4").Value = 0.2136872993148446
$ws.Range("O4").Value = 0.2136744744392
$ws.Range("V4").Value = 955500195836
$ws.Range("W4").Value = 955751087452

# Row 5 (LICs)
$ws.Range("D5").Value = 2.973125562628548
$ws.Range("E5").Value = 4.34634358097193
$ws.Range("F5").Value = 14.86562781314273
$ws.Range("G5").Value = 21.73171790485965

# Row 6 (LMICs)
$ws.Range("D6").Value = 1.387800792749046
$ws.Range("E6").Value = 3.124031414955012
$ws.Range("F6").Value = 6.939003963745232
$ws.Range("G6").Value = 15.62015707477506
$ws.Range("N6").Value = 2.695101883716374
$ws.Range("O6").Value = 0.4751809436654225
$ws.Range("V6").Value = 1074779159
$ws.Range("W6").Value = 16007585528

# Row 7 (Tax haven)
$ws.Range("D7").Value = 5.508305901430718
$ws.Range("E7").Value = 4.152435081943773
$ws.Range("F7").Value = 27.54152950715359
$ws.Range("G7").Value = 20.76217540971887
$ws.Range("N7").Value = -0.4109354431633722
$ws.Range("O7").Value = -0.3565134145334812
$ws.Range("V7").Value = -6943500091
$ws.Range("W7").Value = -7824629507

# Row 8 (UMICs)
$ws.Range("D8").Value = 0.2682359336613213
$ws.Range("E8").Value = 0.5695999317677495
$ws.Range("F8").Value = 1.341179668306607
$ws.Range("G8").Value = 2.847999658838746
$ws.Range("N8").Value = 0.4626563873828087
$ws.Range("O8").Value = 0.4202158248625423
$ws.Range("V8").Value = 37492196667
$ws.Range("W8").Value = 43740019703
